$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update store-related property values on Sheet1.
$ws1.Range("B4").Value = 2
$ws1.Range("B5").Value = "KA63502R01"
$ws1.Range("B6").Value = 63502
$ws1.Range("B8").Value = "xx.xx.xxx.xxx"
$ws1.Range("B9").Value = "KA63502R02"

# Recalculate so Sheet2's formula-derived cached values refresh too.
$excel.CalculateFullRebuild()
